$d = $word.ActiveDocument

# 1. Update the date in the header paragraph
$d.Content.Find.Execute("03.06.24", $true, $false, $false, $false, $false,
                         $true, 1, $false, "02.06.24", 2)

# 2. Replace the paper title (paragraph 2)
$d.Paragraphs(2).Range.Text = "LLaMA-NAS: Efficient Neural Architecture Search for Large Language Models"

# 3. Replace the review body text (paragraph 3)
$d.Paragraphs(3).Range.Text = "פעם הנושא של Neural Architecture Search או NAS בקצרה שעסק בחיפוש לאחר ארכיטקטורה אופטימלית של רשת נוירונים עבור משימה/משימות/דומיין היה די פופולרי אך בשנים האחרונות התחום נמצא בדעיכה. אני שמח שנתקלתי במאמר הזה שמנסה לפתח NAS עבור מודלי שפה. אני זוכר מאמרים די מגניבים שמשתמשים בשיטות RL די מגניבות לכך. אולי בעתיד NAS תהפוך למתחרה רציניות של שיטות פרונינג וקוונטיזציה. "

# 4. Split paragraph 4 ("מאמר: ...") into two paragraphs:
#    the first now holds a lone combining RAFE mark, the second the new link
$p4 = $d.Paragraphs(4)
$p4.Range.Text = "ֿ"
$p4.Range.InsertParagraphAfter()
$d.Paragraphs(5).Range.Text = "מאמר: https://arxiv.org/abs/2405.18377"
